# Added notes extraction feature for PPTX
#
# Slide 1 previously had no speaker notes at all. Add a notes page to it
# with multi-paragraph speaker-notes text (the "notes extraction" sample
# content), mirroring the author's edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$notesPage = $s.NotesPage
$notesBody = $notesPage.Shapes.AddPlaceholder(2)

$notesBody.TextFrame.TextRange.Text = "Testing Multiline Notes.`nTo be extracted here.`nMultiline notes extracted."
